$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bug fix: the T (vx correction) and U (vy correction) columns were
# incorrectly referencing the left-hand simulation block's E/F columns
# instead of their own right-hand block's R/S columns. Fix the formulas
# for rows 4:25 and let Excel recompute relative references per-row.
$ws.Range("T4:T25").FormulaR1C1 = "=RC[-2]+(gx-op*RC[-2]+wx)*dt/2"
$ws.Range("U4:U25").FormulaR1C1 = "=RC[-2]+(gy-op*RC[-2]+wy)*dt/2"

# Restore the active selection as left by the author after the fix.
$ws.Range("AA21").Select()
